$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update time-slot values in column C (rows 6-11).
# C6 stays the same; C7 and C8 get new values; C9-C11 stay the same
# (but because the underlying shared string for the removed slot
# shifts, Excel will simply show the new literal text here).
$ws.Range("C7").Value = "12:50-12:55"
$ws.Range("C8").Value = "12:55-13:0"

# Move the active selection to C16, matching the saved view state.
$ws.Range("C16").Select()
